$wb = $excel.ActiveWorkbook

# ---- ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H80").Value = 1198.7
$ws.Range("J80").Value = 1275
$ws.Range("L80").Value = 3825
$ws.Range("N80").Value = -5821
$ws.Range("H83").Value = 1198.7
$ws.Range("J83").Value = 1275
$ws.Range("L83").Value = 11475
$ws.Range("N83").Value = -21459
$ws.Range("H137").Value = 1432.6666
$ws.Range("I137").Value = 1459.2667
$ws.Range("J137").Value = 1299.6666
$ws.Range("K137").Value = 4377.800099999999
$ws.Range("L137").Value = 3898.9998
$ws.Range("M137").Value = -1827.800099999999
$ws.Range("N137").Value = -8998.9998

# ---- ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 3961.2407
$ws.Range("I32").Value = 3610.157
$ws.Range("J32").Value = 9929.666999999999
$ws.Range("K32").Value = 3610.157
$ws.Range("L32").Value = 9929.666999999999
$ws.Range("M32").Value = -3323.157
$ws.Range("N32").Value = -10503.667
$ws.Range("H74").Value = 6145.5264
$ws.Range("I74").Value = 1401.129
$ws.Range("J74").Value = 27156.428
$ws.Range("K74").Value = 1401.129
$ws.Range("L74").Value = 27156.428
$ws.Range("M74").Value = -527.1289999999999
$ws.Range("N74").Value = -28904.428
$ws.Range("H77").Value = 6145.5264
$ws.Range("I77").Value = 1401.129
$ws.Range("J77").Value = 27156.428
$ws.Range("K77").Value = 7005.645
$ws.Range("L77").Value = 135782.14
$ws.Range("M77").Value = -2637.645
$ws.Range("N77").Value = -144518.14
$ws.Range("H97").Value = 1057.25
$ws.Range("I97").Value = 811.3043
$ws.Range("K97").Value = 811.3043
$ws.Range("M97").Value = -315.3043
$ws.Range("H122").Value = 2382.6072
$ws.Range("J122").Value = 2850
$ws.Range("L122").Value = 8550
$ws.Range("N122").Value = -13450
$ws.Range("H132").Value = 2528.843
$ws.Range("I132").Value = 2427.2327
$ws.Range("K132").Value = 7281.6981
$ws.Range("M132").Value = -4751.6981

# ---- BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 2701.0286
$ws.Range("I20").Value = 2087.3333
$ws.Range("J20").Value = 3621.5715
$ws.Range("K20").Value = 2087.3333
$ws.Range("L20").Value = 3621.5715
$ws.Range("M20").Value = -1840.3333
$ws.Range("N20").Value = -4115.5715
$ws.Range("H94").Value = 1373.3
$ws.Range("I94").Value = 1498
$ws.Range("J94").Value = 874.5
$ws.Range("K94").Value = 1498
$ws.Range("L94").Value = 874.5
$ws.Range("M94").Value = -1047
$ws.Range("N94").Value = -1776.5

# ---- CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H4").Value = 350972.22
$ws.Range("J4").Value = 756000
$ws.Range("L4").Value = 756000
$ws.Range("N4").Value = -756224
$ws.Range("H22").Value = 670.5714
$ws.Range("I22").Value = 670.5714
$ws.Range("K22").Value = 670.5714
$ws.Range("M22").Value = -320.5714
$ws.Range("H132").Value = 4110.625
$ws.Range("I132").Value = 3942.08
$ws.Range("K132").Value = 11826.24
$ws.Range("M132").Value = -9296.24

# ---- CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 54.464287
$ws.Range("I2").Value = 56.2
$ws.Range("K2").Value = 337.2
$ws.Range("M2").Value = -224.2
$ws.Range("H17").Value = 39.6
$ws.Range("J17").Value = 54.5
$ws.Range("L17").Value = 163.5
$ws.Range("N17").Value = -501.5
$ws.Range("H38").Value = 81.888885
$ws.Range("I38").Value = 83
$ws.Range("K38").Value = 249
$ws.Range("M38").Value = 98
$ws.Range("H68").Value = 12502500
$ws.Range("I68").Value = 2000
$ws.Range("K68").Value = 6000
$ws.Range("M68").Value = -5189
$ws.Range("H71").Value = 12502500
$ws.Range("I71").Value = 2000
$ws.Range("K71").Value = 18000
$ws.Range("M71").Value = -13944
$ws.Range("H80").Value = 0
$ws.Range("J80").Value = 0
$ws.Range("L80").ClearContents()
$ws.Range("N80").Value = 0
$ws.Range("H83").Value = 0
$ws.Range("J83").Value = 0
$ws.Range("L83").ClearContents()
$ws.Range("N83").Value = 0
$ws.Range("H134").Value = 6838.4
$ws.Range("I134").Value = 5589.091
$ws.Range("K134").Value = 16767.273
$ws.Range("M134").Value = -11697.273

# ---- GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H107").Value = 412.72726
$ws.Range("I107").Value = 365.2
$ws.Range("K107").Value = 365.2
$ws.Range("M107").Value = 1554.8
$ws.Range("H132").Value = 2845.1072
$ws.Range("J132").Value = 3725.6
$ws.Range("M132").Value = -16236.8
$ws.Range("N132").Value = -16236.8

# ---- LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 1469
$ws.Range("I22").Value = 1990
$ws.Range("J22").Value = 1411.1111
$ws.Range("K22").Value = 1990
$ws.Range("L22").Value = 1411.1111
$ws.Range("M22").Value = -1695
$ws.Range("N22").Value = -2001.1111
$ws.Range("H27").Value = 1469
$ws.Range("I27").Value = 1990
$ws.Range("J27").Value = 1411.1111
$ws.Range("K27").Value = 1990
$ws.Range("L27").Value = 1411.1111
$ws.Range("M27").Value = -1883
$ws.Range("N27").Value = -1625.1111
$ws.Range("H42").Value = 17000
$ws.Range("I42").Value = 16800
$ws.Range("J42").Value = 19000
$ws.Range("K42").Value = 16800
$ws.Range("L42").Value = 19000
$ws.Range("M42").Value = -16237
$ws.Range("N42").Value = -20126
$ws.Range("H46").Value = 2536
$ws.Range("I46").Value = 1993.3334
$ws.Range("J46").Value = 3350
$ws.Range("K46").Value = 1993.3334
$ws.Range("L46").Value = 3350
$ws.Range("M46").Value = -1805.3334
$ws.Range("N46").Value = -3726
$ws.Range("H49").Value = 17000
$ws.Range("I49").Value = 16800
$ws.Range("J49").Value = 19000
$ws.Range("K49").Value = 16800
$ws.Range("L49").Value = 19000
$ws.Range("M49").Value = -16653
$ws.Range("N49").Value = -19294
$ws.Range("H122").Value = 5085.5
$ws.Range("I122").Value = 4352.2856
$ws.Range("J122").Value = 7285.143
$ws.Range("K122").Value = 13056.8568
$ws.Range("L122").Value = 21855.429
$ws.Range("M122").Value = -10606.8568
$ws.Range("N122").Value = -26755.429
$ws.Range("H132").Value = 5857.7144
$ws.Range("I132").Value = 5251
$ws.Range("J132").Value = 6666.6665
$ws.Range("K132").Value = 15753
$ws.Range("L132").Value = 19999.9995
$ws.Range("M132").Value = -13223
$ws.Range("N132").Value = -25059.9995
$ws.Range("H133").Value = 59579.8
$ws.Range("J133").Value = 59579.8
$ws.Range("L133").Value = 59579.8
$ws.Range("N133").Value = -64639.8

# ---- WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H5").Value = 999
$ws.Range("J5").Value = 999
$ws.Range("L5").Value = 999
$ws.Range("N5").Value = -1223
$ws.Range("H14").Value = 370729.38
$ws.Range("I14").Value = 463093.25
$ws.Range("J14").Value = 1273.8334
$ws.Range("K14").Value = 463093.25
$ws.Range("L14").Value = 1273.8334
$ws.Range("M14").Value = -462925.25
$ws.Range("N14").Value = -1609.8334
$ws.Range("H132").Value = 1659.9762
$ws.Range("I132").Value = 1592.975
$ws.Range("J132").Value = 3000
$ws.Range("K132").Value = 4778.924999999999
$ws.Range("L132").Value = 9000
$ws.Range("M132").Value = -2248.924999999999
$ws.Range("N132").Value = -14060
